$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (col C, D, F, H changed; values chosen so stored
# OOXML width equals target after Excel's internal +0.8333 padding)
$ws.Columns.Item(3).ColumnWidth = 81 - 0.8333333333333333   # C: 69 -> 81
$ws.Columns.Item(4).ColumnWidth = 50 - 0.8333333333333333   # D: 30 -> 50
$ws.Columns.Item(6).ColumnWidth = 17 - 0.8333333333333333   # F: 16 -> 17
$ws.Columns.Item(8).ColumnWidth = 29 - 0.8333333333333333   # H: 41 -> 29

# Column A holds opportunity IDs stored as text (not numbers) in the source data,
# so pre-format the data rows as Text before writing the numeric-looking IDs.
$ws.Range("A2:A18").NumberFormat = "@"

$data = @(
    @('1329542', 'https://aiesec.org/opportunity/global-talent/1329542', 'Junior Software Developer Threesixtymedia', '52 Aachen, Germany', 'No', '6 applicants', '3 - 6 Months', 'Threesixtymedia'),
    @('1329526', 'https://aiesec.org/opportunity/global-talent/1329526', '[Impact Brazil] - Supplier Development Procurement Intern', 'Joinville - Pirabeiraba, Joinville - SC, Brasil', 'No', '1 applicant', '3 - 6 Months', 'Nidec Global Appliance'),
    @('1329520', 'https://aiesec.org/opportunity/global-talent/1329520', 'Sales Officer', 'Dehiwala-Mount Lavinia, Sri Lanka', 'No', '0 applicants', '3 - 6 Months', 'Pedro Barn pvt ltd'),
    @('1329480', 'https://aiesec.org/opportunity/global-talent/1329480', '3D Modeling intern', 'Amman, Jordan', 'No', '1 applicant', '9 - 12 Weeks', 'MKHAMAWI EST'),
    @('1329458', 'https://aiesec.org/opportunity/global-talent/1329458', 'Digital Marketing Intern', 'Visakhapatnam, Andhra Pradesh, India', 'No', '3 applicants', '9 - 12 Weeks', 'Media3 Training'),
    @('1329457', 'https://aiesec.org/opportunity/global-talent/1329457', 'BA Intern', 'Visakhapatnam, Andhra Pradesh, India', 'No', '1 applicant', '9 - 12 Weeks', 'Media3 Training'),
    @('1329177', 'https://aiesec.org/opportunity/global-talent/1329177', '[Impact Porto Alegre] - Social Media Intern – Intercultural Engagement Program', 'Porto Alegre, RS, Brazil', 'No', '16 applicants', '9 - 12 Weeks', 'ESCOLA GIORDANO BRUNO LTDA'),
    @('1328845', 'https://aiesec.org/opportunity/global-talent/1328845', 'Business Development Intern', 'Hyderabad, Telangana, India', 'No', '8 applicants', '3 - 6 Months', 'Vigilare biopharma Pvt Ltd'),
    @('1327949', 'https://aiesec.org/opportunity/global-talent/1327949', 'Taste Hungary | Football Data Analyst', 'Budapest, Magyarország', 'No', '79 applicants', '9 - 12 Weeks', 'ACE Advisory'),
    @('1327439', 'https://aiesec.org/opportunity/global-talent/1327439', 'Assistant', 'Tunis, Tunisie', 'No', '33 applicants', '9 - 12 Weeks', 'SSK events'),
    @('1326906', 'https://aiesec.org/opportunity/global-talent/1326906', 'Sales Manager', 'Mersin, Akdeniz/Mersin, Türkiye', 'No', '72 applicants', '3 - 6 Months', 'İlke Sports'),
    @('1326741', 'https://aiesec.org/opportunity/global-talent/1326741', 'Business Development Intern', 'Hyderabad, Telangana, India', 'No', '7 applicants', '9 - 12 Weeks', 'Vigilare biopharma Pvt Ltd'),
    @('1326373', 'https://aiesec.org/opportunity/global-talent/1326373', 'Global Mobility Trainee', 'San Francisco, CA, USA', 'No', '245 applicants', '6 - 18 Months', 'Vialto Partners US'),
    @('1325604', 'https://aiesec.org/opportunity/global-talent/1325604', 'International Business & Innovation Analyst', '4520 Santa Maria da Feira, Portugal', 'No', '163 applicants', '9 - 12 Weeks', 'M2K Consultoria'),
    @('1321054', 'https://aiesec.org/opportunity/global-talent/1321054', 'International Sales Representetive Russian Speaker', 'Maslak, Sarıyer/İstanbul, Türkiye', 'No', '15 applicants', '6 - 18 Months', 'Esvita Clinic'),
    @('1321053', 'https://aiesec.org/opportunity/global-talent/1321053', 'International Sales Representetive German Speaker', 'Maslak, Sarıyer/İstanbul, Türkiye', 'No', '17 applicants', '6 - 18 Months', 'Esvita Clinic'),
    @('1315274', 'https://aiesec.org/opportunity/global-talent/1315274', '[Impact Florianópolis]- Social Media', 'São Miguel do Oeste, SC, 89900-000, Brasil', 'No', '27 applicants', '9 - 12 Weeks', 'KNN Idiomas'),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

Write-Host "Updated rows 2-18"
